# feat: add 2022-Q1 data
#
# - Repurpose the existing "总计" sheet (sheetId 3 / rId3) to become the new
#   "2022-Q1" per-fund holdings sheet, keeping its sheetId/rId (matches the
#   target OOXML, which reuses sheetId="3" / r:id="rId3" for "2022-Q1").
# - Add a brand new "总计" sheet right after it (gets the next sheetId / rId,
#   i.e. sheetId="4" / r:id="rId4"), and rebuild the summary table with the
#   new 2022-Q1 row inserted at the top.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Grab a reference to the formatted header style (bold + border +
#    centered) that's already used on the "2021-Q4" sheet's header row /
#    index column, so newly-added cells match the sheet's existing look.
# ---------------------------------------------------------------------
$styleSrc = $wb.Worksheets.Item("2021-Q4").Range("B1")

function Set-HeaderStyle($range) {
    $styleSrc.Copy()
    $range.PasteSpecial(-4122)  # xlPasteFormats
}

function Set-TextValue($cell, $text) {
    # Force text storage (matches the source data, which stores numeric-
    # looking strings like "5.78" as text, not numbers).
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# ---------------------------------------------------------------------
# 2. Repurpose the old "总计" sheet -> "2022-Q1" fund holdings sheet.
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"
$q1.Cells.Clear()

# Header row
Set-HeaderStyle($q1.Range("B1:H1"))
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Fund rows: code, name, scale, position, ratio, value, rank
$q1Rows = @(
    @("006440", "中信建投中证500指数增强A", "5.78", "94.71", "1.11", "0.0642", 4),
    @("006441", "中信建投中证500指数增强C", "3.11", "94.71", "1.11", "0.0345", 4),
    @("006227", "华宝科技先锋混合A",         "0.46", "91.73", "2.66", "0.0122", 8),
    @("002802", "广发东财大数据精选灵活配置混合", "0.41", "55.13", "2.28", "0.0093", 6),
    @("010842", "华宝科技先锋混合C",         "0.04", "91.73", "2.66", "0.0011", 8)
)

$r = 2
foreach ($row in $q1Rows) {
    Set-HeaderStyle($q1.Range("A$r"))
    $q1.Range("A$r").Value = ($r - 2)
    Set-TextValue $q1.Range("B$r") $row[0]
    $q1.Range("C$r").Value = $row[1]
    Set-TextValue $q1.Range("D$r") $row[2]
    Set-TextValue $q1.Range("E$r") $row[3]
    Set-TextValue $q1.Range("F$r") $row[4]
    Set-TextValue $q1.Range("G$r") $row[5]
    $q1.Range("H$r").Value = $row[6]
    $r = $r + 1
}

# ---------------------------------------------------------------------
# 3. Insert a fresh "总计" sheet right after "2022-Q1".
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

Set-HeaderStyle($total.Range("B1:D1"))
$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$totalRows = @(
    @("2022-Q1", 5, 0.12),
    @("2021-Q4", 1, 0.16),
    @("2021-Q3", 4, 0.53)
)

$r = 2
foreach ($row in $totalRows) {
    Set-HeaderStyle($total.Range("A$r"))
    $total.Range("A$r").Value = ($r - 2)
    $total.Range("B$r").Value = $row[0]
    $total.Range("C$r").Value = $row[1]
    $total.Range("D$r").Value = $row[2]
    $r = $r + 1
}

$wb.Worksheets.Item("2021-Q3").Activate()
